# Add a new slide at the end of the deck, using the "Title and Content"
# layout (the same layout already used by the existing diagram slide),
# and fill in the body placeholder with the new text.

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)   # 2 = ppLayoutText ("Title and Content")

# Shape 1 is the Title placeholder - leave it blank.
# Shape 2 is the body/content placeholder - set its text.
$s.Shapes.Item(2).TextFrame.TextRange.Text = "For test"
